$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"; $ws.Range("D2").Value = "285.28"
$ws.Range("E2").NumberFormat = "@"; $ws.Range("E2").Value = "2.30%"
$ws.Range("E3").NumberFormat = "@"; $ws.Range("E3").Value = "4.75%"
$ws.Range("E4").NumberFormat = "@"; $ws.Range("E4").Value = "5.71%"
$ws.Range("D5").NumberFormat = "@"; $ws.Range("D5").Value = "0.06683"
$ws.Range("E5").NumberFormat = "@"; $ws.Range("E5").Value = "4.75%"
$ws.Range("D6").NumberFormat = "@"; $ws.Range("D6").Value = "7.325"
$ws.Range("E6").NumberFormat = "@"; $ws.Range("E6").Value = "4.13%"
$ws.Range("D7").NumberFormat = "@"; $ws.Range("D7").Value = "3.385"
$ws.Range("E7").NumberFormat = "@"; $ws.Range("E7").Value = "2.51%"
$ws.Range("D8").NumberFormat = "@"; $ws.Range("D8").Value = "1.351"
$ws.Range("E8").NumberFormat = "@"; $ws.Range("E8").Value = "4.28%"
$ws.Range("D9").NumberFormat = "@"; $ws.Range("D9").Value = "0.9360"
$ws.Range("E9").NumberFormat = "@"; $ws.Range("E9").Value = "4.98%"
$ws.Range("D10").NumberFormat = "@"; $ws.Range("D10").Value = "0.1572"
$ws.Range("E10").NumberFormat = "@"; $ws.Range("E10").Value = "3.03%"
$ws.Range("D11").NumberFormat = "@"; $ws.Range("D11").Value = "0.06524"
$ws.Range("E11").NumberFormat = "@"; $ws.Range("E11").Value = "15.77%"
$ws.Range("D12").NumberFormat = "@"; $ws.Range("D12").Value = "0.07682"
$ws.Range("E12").NumberFormat = "@"; $ws.Range("E12").Value = "2.45%"
$ws.Range("D13").NumberFormat = "@"; $ws.Range("D13").Value = "0.02888"
$ws.Range("E13").NumberFormat = "@"; $ws.Range("E13").Value = "-0.86%"
$ws.Range("D14").NumberFormat = "@"; $ws.Range("D14").Value = "0.08976"
$ws.Range("E14").NumberFormat = "@"; $ws.Range("E14").Value = "-0.04%"
$ws.Range("D15").NumberFormat = "@"; $ws.Range("D15").Value = "0.001591"
$ws.Range("E15").NumberFormat = "@"; $ws.Range("E15").Value = "0.95%"
$ws.Range("D16").NumberFormat = "@"; $ws.Range("D16").Value = "0.04477"
$ws.Range("E16").NumberFormat = "@"; $ws.Range("E16").Value = "2.22%"
$ws.Range("D17").NumberFormat = "@"; $ws.Range("D17").Value = "0.0006446"
$ws.Range("E17").NumberFormat = "@"; $ws.Range("E17").Value = "0.91%"
$ws.Range("D18").NumberFormat = "@"; $ws.Range("D18").Value = "0.006544"
$ws.Range("E18").NumberFormat = "@"; $ws.Range("E18").Value = "7.38%"
$ws.Range("E19").NumberFormat = "@"; $ws.Range("E19").Value = "0.24%"
$ws.Range("D20").NumberFormat = "@"; $ws.Range("D20").Value = "2.235"
$ws.Range("E20").NumberFormat = "@"; $ws.Range("E20").Value = "-2.28%"
$ws.Range("D21").NumberFormat = "@"; $ws.Range("D21").Value = "0.3206"
$ws.Range("E21").NumberFormat = "@"; $ws.Range("E21").Value = "1.87%"
$ws.Range("D22").NumberFormat = "@"; $ws.Range("D22").Value = "0.1306"
$ws.Range("E22").NumberFormat = "@"; $ws.Range("E22").Value = "-3.25%"
$ws.Range("D23").NumberFormat = "@"; $ws.Range("D23").Value = "4.050"
$ws.Range("E23").NumberFormat = "@"; $ws.Range("E23").Value = "3.73%"
$ws.Range("E24").NumberFormat = "@"; $ws.Range("E24").Value = "1.17%"
$ws.Range("D25").NumberFormat = "@"; $ws.Range("D25").Value = "0.001178"
$ws.Range("E25").NumberFormat = "@"; $ws.Range("E25").Value = "0.29%"
$ws.Range("D26").NumberFormat = "@"; $ws.Range("D26").Value = "0.004470"
$ws.Range("E26").NumberFormat = "@"; $ws.Range("E26").Value = "4.41%"
$ws.Range("D27").NumberFormat = "@"; $ws.Range("D27").Value = "0.0001244"
$ws.Range("E27").NumberFormat = "@"; $ws.Range("E27").Value = "5.69%"
$ws.Range("E28").NumberFormat = "@"; $ws.Range("E28").Value = "-2.35%"
$ws.Range("D40").NumberFormat = "@"; $ws.Range("D40").Value = "0.04194"
$ws.Range("E40").NumberFormat = "@"; $ws.Range("E40").Value = "3.86%"
$ws.Range("D41").NumberFormat = "@"; $ws.Range("D41").Value = "0.006698"
$ws.Range("E41").NumberFormat = "@"; $ws.Range("E41").Value = "-0.31%"
$ws.Range("D42").NumberFormat = "@"; $ws.Range("D42").Value = "0.1247"
$ws.Range("E42").NumberFormat = "@"; $ws.Range("E42").Value = "-11.29%"
$ws.Range("D43").NumberFormat = "@"; $ws.Range("D43").Value = "0.002011"
$ws.Range("E43").NumberFormat = "@"; $ws.Range("E43").Value = "-1.69%"
$ws.Range("D44").NumberFormat = "@"; $ws.Range("D44").Value = "0.01214"
$ws.Range("E44").NumberFormat = "@"; $ws.Range("E44").Value = "8.87%"
$ws.Range("D45").NumberFormat = "@"; $ws.Range("D45").Value = "0.00005667"
$ws.Range("E45").NumberFormat = "@"; $ws.Range("E45").Value = "2.24%"
$ws.Range("E46").NumberFormat = "@"; $ws.Range("E46").Value = "25.93%"
$ws.Range("E47").NumberFormat = "@"; $ws.Range("E47").Value = "-29.47%"
